$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Definition (column C) and Parent (column D) cells for BCIO population
# statistic rows: switch from the old "The aggregate of X in a population." /
# "data item" parent pattern to "A population statistic about X." with a more
# specific parent term.
$data = @(
    @{ Row = 234; C = 'A population statistic about ability to comprehend spoken intervention language.'; D = 'linguistic capability  population statistic' },
    @{ Row = 237; C = 'A population statistic about ability to read in intervention language.'; D = 'linguistic capability  population statistic' },
    @{ Row = 240; C = 'A population statistic about ability to speak in intervention language.'; D = 'linguistic capability  population statistic' },
    @{ Row = 243; C = 'A population statistic about ability to write in intervention language.'; D = 'linguistic capability  population statistic' },
    @{ Row = 246; C = 'A population statistic about achieved bachelor''s degree or equivalent level.'; D = 'highest level of formal educational qualification achieved population statistic' },
    @{ Row = 249; C = 'A population statistic about achieved doctoral or equivalent level education.'; D = 'highest level of formal educational qualification achieved population statistic' },
    @{ Row = 252; C = 'A population statistic about achieved early childhood education.'; D = 'highest level of formal educational qualification achieved population statistic' },
    @{ Row = 255; C = 'A population statistic about achieved lower secondary education.'; D = 'highest level of formal educational qualification achieved population statistic' },
    @{ Row = 258; C = 'A population statistic about achieved master''s or equivalent level.'; D = 'highest level of formal educational qualification achieved population statistic' },
    @{ Row = 261; C = 'A population statistic about achieved primary education.'; D = 'highest level of formal educational qualification achieved population statistic' },
    @{ Row = 264; C = 'A population statistic about achieved upper secondary education .'; D = 'population statistic' },
    @{ Row = 267; C = 'A population statistic about adoptive brother.'; D = 'adoptive sibling population statistic' },
    @{ Row = 270; C = 'A population statistic about adoptive child.'; D = 'child relation population statistic' },
    @{ Row = 273; C = 'A population statistic about adoptive daughter.'; D = 'adoptive child population statistic' },
    @{ Row = 276; C = 'A population statistic about adoptive father.'; D = 'adoptive parent population statistic' },
    @{ Row = 279; C = 'A population statistic about adoptive mother.'; D = 'adoptive parent population statistic' },
    @{ Row = 282; C = 'A population statistic about adoptive parent.'; D = 'parent population statistic' },
    @{ Row = 285; C = 'A population statistic about adoptive sibling.'; D = 'sibling population statistic' },
    @{ Row = 288; C = 'A population statistic about adoptive sister.'; D = 'adoptive sibling population statistic' },
    @{ Row = 291; C = 'A population statistic about adoptive son.'; D = 'adoptive child population statistic' },
    @{ Row = 294; C = 'A population statistic about adult.'; D = 'person population statistic' },
    @{ Row = 297; C = 'A population statistic about agreed rent-free occupier.'; D = 'rent-free occupier population statistic' },
    @{ Row = 300; C = 'A population statistic about asexual.'; D = 'sexual orientation population statistic' },
    @{ Row = 303; C = 'A population statistic about aunt.'; D = 'family member population statistic' },
    @{ Row = 306; C = 'A population statistic about biological brother.'; D = 'biological sibling population statistic' },
    @{ Row = 309; C = 'A population statistic about biological child.'; D = 'child relation population statistic' },
    @{ Row = 312; C = 'A population statistic about biological daughter.'; D = 'biological child population statistic' },
    @{ Row = 315; C = 'A population statistic about biological father.'; D = 'biological parent population statistic' },
    @{ Row = 318; C = 'A population statistic about biological mother.'; D = 'biological parent population statistic' },
    @{ Row = 321; C = 'A population statistic about biological parent.'; D = 'parent population statistic' },
    @{ Row = 324; C = 'A population statistic about biological sex.'; D = 'bodily quality population statistic' },
    @{ Row = 327; C = 'A population statistic about biological sibling.'; D = 'sibling population statistic' },
    @{ Row = 330; C = 'A population statistic about biological sister.'; D = 'biological sibling population statistic' },
    @{ Row = 333; C = 'A population statistic about biological son.'; D = 'biological child population statistic' },
    @{ Row = 336; C = 'A population statistic about bisexual.'; D = 'sexual orientation population statistic' },
    @{ Row = 339; C = 'A population statistic about caregiving role.'; D = 'role population statistic' },
    @{ Row = 342; C = 'A population statistic about caste membership.'; D = 'personal attribute population statistic' },
    @{ Row = 345; C = 'A population statistic about child.'; D = 'person population statistic' },
    @{ Row = 348; C = 'A population statistic about child relation.'; D = 'family member population statistic' },
    @{ Row = 351; C = 'A population statistic about cisgender.'; D = 'gender identity population statistic' },
    @{ Row = 354; C = 'A population statistic about country of birth.'; D = 'geographic location population statistic' },
    @{ Row = 357; C = 'A population statistic about cousin.'; D = 'family member population statistic' },
    @{ Row = 360; C = 'A population statistic about disabled.'; D = 'personal attribute population statistic' },
    @{ Row = 363; C = 'A population statistic about discipline of current programme of study or training.'; D = 'expertise discipline population statistic' },
    @{ Row = 366; C = 'A population statistic about discipline of highest level of formal educational qualification achieved.'; D = 'expertise discipline population statistic' },
    @{ Row = 369; C = 'A population statistic about divorced or separated.'; D = 'relationship status population statistic' },
    @{ Row = 372; C = 'A population statistic about doctoral student role.'; D = 'higher education student role population statistic' },
    @{ Row = 375; C = 'A population statistic about employed.'; D = 'employment status population statistic' },
    @{ Row = 378; C = 'A population statistic about employed full time.'; D = 'employment status population statistic' },
    @{ Row = 381; C = 'A population statistic about employed in shift work.'; D = 'employment status population statistic' },
    @{ Row = 384; C = 'A population statistic about employed part time.'; D = 'employment status population statistic' },
    @{ Row = 387; C = 'A population statistic about employment status.'; D = 'quality population statistic' },
    @{ Row = 390; C = 'A population statistic about ethnic group membership.'; D = 'self-identity population statistic' },
    @{ Row = 393; C = 'A population statistic about expertise discipline.'; D = 'specifically dependent continuant population statistic' },
    @{ Row = 396; C = 'A population statistic about family member.'; D = 'person population statistic' },
    @{ Row = 399; C = 'A population statistic about father.'; D = 'parent population statistic' },
    @{ Row = 402; C = 'A population statistic about female biological sex.'; D = 'biological sex population statistic' },
    @{ Row = 405; C = 'A population statistic about female gender.'; D = 'gender identity population statistic' },
    @{ Row = 408; C = 'A population statistic about foster brother.'; D = 'foster sibling population statistic' },
    @{ Row = 411; C = 'A population statistic about foster child.'; D = 'child relation population statistic' },
    @{ Row = 414; C = 'A population statistic about foster daughter.'; D = 'foster child population statistic' },
    @{ Row = 417; C = 'A population statistic about foster father.'; D = 'foster parent population statistic' },
    @{ Row = 420; C = 'A population statistic about foster mother.'; D = 'foster parent population statistic' },
    @{ Row = 423; C = 'A population statistic about foster parent.'; D = 'parent population statistic' },
    @{ Row = 426; C = 'A population statistic about foster sibling.'; D = 'sibling population statistic' },
    @{ Row = 429; C = 'A population statistic about foster sister.'; D = 'foster sibling population statistic' },
    @{ Row = 432; C = 'A population statistic about foster son.'; D = 'foster child population statistic' },
    @{ Row = 435; C = 'A population statistic about gender identity.'; D = 'self-identity population statistic' },
    @{ Row = 438; C = 'A population statistic about graduate student role.'; D = 'higher education student role population statistic' },
    @{ Row = 441; C = 'A population statistic about grandfather.'; D = 'grandparent population statistic' },
    @{ Row = 444; C = 'A population statistic about grandmother.'; D = 'grandparent population statistic' },
    @{ Row = 447; C = 'A population statistic about grandparent.'; D = 'family member population statistic' },
    @{ Row = 450; C = 'A population statistic about having enacted a behaviour.'; D = 'personal history part population statistic' },
    @{ Row = 453; C = 'A population statistic about health insurance policy holder role.'; D = 'policy holder role population statistic' },
    @{ Row = 456; C = 'A population statistic about health status attribute.'; D = 'personal attribute population statistic' },
    @{ Row = 459; C = 'A population statistic about heterosexual.'; D = 'sexual orientation population statistic' },
    @{ Row = 462; C = 'A population statistic about higher education student role.'; D = 'student or trainee role population statistic' },
    @{ Row = 465; C = 'A population statistic about highest level of formal educational qualification achieved.'; D = 'personal attribute population statistic' },
    @{ Row = 468; C = 'A population statistic about history of exposure to an occupational hazard.'; D = 'personal history part population statistic' },
    @{ Row = 471; C = 'A population statistic about history of exposure to childhood maltreatment.'; D = 'personal history part population statistic' },
    @{ Row = 474; C = 'A population statistic about homeless person.'; D = 'person population statistic' },
    @{ Row = 477; C = 'A population statistic about homemaker status.'; D = 'personal attribute population statistic' },
    @{ Row = 480; C = 'A population statistic about homosexual.'; D = 'sexual orientation population statistic' },
    @{ Row = 483; C = 'A population statistic about household income.'; D = 'object aggregate population statistic' },
    @{ Row = 490; C = 'A population statistic about human age.'; D = 'personal attribute population statistic' },
    @{ Row = 495; C = 'A population statistic about immigrant.'; D = 'person population statistic' },
    @{ Row = 498; C = 'A population statistic about in a legal marriage or union.'; D = 'relationship status population statistic' },
    @{ Row = 501; C = 'A population statistic about in a stable or common law relationship.'; D = 'relationship status population statistic' },
    @{ Row = 504; C = 'A population statistic about in permanent employment.'; D = 'employment status population statistic' },
    @{ Row = 507; C = 'A population statistic about in short term or temporary employment with known conditions.'; D = 'employment status population statistic' },
    @{ Row = 510; C = 'A population statistic about in uncertain employment.'; D = 'employment status population statistic' },
    @{ Row = 513; C = 'A population statistic about income-related welfare benefit.'; D = 'individual income population statistic' },
    @{ Row = 516; C = 'A population statistic about independently wealthy status.'; D = 'personal attribute population statistic' },
    @{ Row = 519; C = 'A population statistic about individual human behaviour.'; D = 'bodily process population statistic' },
    @{ Row = 526; C = 'A population statistic about individual income.'; D = 'personal attribute population statistic' },
    @{ Row = 533; C = 'A population statistic about influencer role.'; D = 'social role population statistic' },
    @{ Row = 536; C = 'A population statistic about informal education student role.'; D = 'student or trainee role population statistic' },
    @{ Row = 539; C = 'A population statistic about inpatient role.'; D = 'patient role population statistic' },
    @{ Row = 542; C = 'A population statistic about insured party role.'; D = 'role population statistic' },
    @{ Row = 545; C = 'A population statistic about interpersonal role.'; D = 'role population statistic' },
    @{ Row = 548; C = 'A population statistic about language proficiency.'; D = 'linguistic capability  population statistic' },
    @{ Row = 555; C = 'A population statistic about linguistic capability.'; D = 'mental capability population statistic' },
    @{ Row = 562; C = 'A population statistic about long-term disabled.'; D = 'disabled population statistic' },
    @{ Row = 565; C = 'A population statistic about male biological sex.'; D = 'biological sex population statistic' },
    @{ Row = 568; C = 'A population statistic about male gender.'; D = 'gender identity population statistic' },
    @{ Row = 571; C = 'A population statistic about masters student role.'; D = 'higher education student role population statistic' },
    @{ Row = 574; C = 'A population statistic about medication use status.'; D = 'health status attribute population statistic' },
    @{ Row = 577; C = 'A population statistic about member of a multi-person household.'; D = 'person population statistic' },
    @{ Row = 580; C = 'A population statistic about member of a multi-person household all related.'; D = 'member of a multi-person household population statistic' },
    @{ Row = 583; C = 'A population statistic about member of a multi-person household not related.'; D = 'member of a multi-person household population statistic' },
    @{ Row = 586; C = 'A population statistic about member of a multi-person household some related.'; D = 'member of a multi-person household population statistic' },
    @{ Row = 589; C = 'A population statistic about member of a multi-person multi-generational household.'; D = 'member of a multi-person household population statistic' },
    @{ Row = 592; C = 'A population statistic about member of a one person household.'; D = 'person population statistic' },
    @{ Row = 595; C = 'A population statistic about mental capability.'; D = 'personal capability population statistic' },
    @{ Row = 602; C = 'A population statistic about mother.'; D = 'parent population statistic' },
    @{ Row = 605; C = 'A population statistic about nephew.'; D = 'family member population statistic' },
    @{ Row = 608; C = 'A population statistic about niece.'; D = 'family member population statistic' },
    @{ Row = 611; C = 'A population statistic about non-gendered identity.'; D = 'self-identity population statistic' },
    @{ Row = 614; C = 'A population statistic about nonbinary gender.'; D = 'gender identity population statistic' },
    @{ Row = 617; C = 'A population statistic about not seeking employment.'; D = 'personal attribute population statistic' },
    @{ Row = 620; C = 'A population statistic about not working for health reasons.'; D = 'personal attribute population statistic' },
    @{ Row = 623; C = 'A population statistic about number of years in education completed.'; D = 'data item population statistic' },
    @{ Row = 628; C = 'A population statistic about occupational role.'; D = 'personal role population statistic' },
    @{ Row = 631; C = 'A population statistic about occupier of employer-provided housing.'; D = 'person population statistic' },
    @{ Row = 634; C = 'A population statistic about organisational role.'; D = 'role population statistic' },
    @{ Row = 637; C = 'A population statistic about other sexual orientation.'; D = 'sexual orientation population statistic' },
    @{ Row = 640; C = 'A population statistic about outpatient role.'; D = 'patient role population statistic' },
    @{ Row = 643; C = 'A population statistic about owner.'; D = 'material entity population statistic' },
    @{ Row = 646; C = 'A population statistic about owner-occupier.'; D = 'person population statistic' },
    @{ Row = 649; C = 'A population statistic about parent.'; D = 'family member population statistic' },
    @{ Row = 652; C = 'A population statistic about parental role.'; D = 'interpersonal role population statistic' },
    @{ Row = 655; C = 'A population statistic about past behaviour .'; D = 'population statistic' },
    @{ Row = 662; C = 'A population statistic about patient role.'; D = 'role population statistic' },
    @{ Row = 665; C = 'A population statistic about personal history of behavioural lapse.'; D = 'personal history part population statistic' },
    @{ Row = 672; C = 'A population statistic about personal history of events that influence behaviour .'; D = 'population statistic' },
    @{ Row = 679; C = 'A population statistic about personal history of intervention exposure for the same outcome.'; D = 'personal history part population statistic' },
    @{ Row = 686; C = 'A population statistic about personal history of intervention exposure for the same outcome behaviour.'; D = 'personal history part population statistic' },
    @{ Row = 693; C = 'A population statistic about personal history of same intervention exposure.'; D = 'personal history part population statistic' },
    @{ Row = 700; C = 'A population statistic about personal history part of intervention exposure.'; D = 'personal history part population statistic' },
    @{ Row = 707; C = 'A population statistic about personal psychological attribute.'; D = 'personal attribute population statistic' },
    @{ Row = 714; C = 'A population statistic about personal vulnerability.'; D = 'disposition population statistic' },
    @{ Row = 721; C = 'A population statistic about personal vulnerability to harmful behaviour.'; D = 'personal vulnerability population statistic' },
    @{ Row = 728; C = 'A population statistic about place of residence.'; D = 'geographic location population statistic' },
    @{ Row = 731; C = 'A population statistic about policy holder role.'; D = 'insured party role population statistic' },
    @{ Row = 734; C = 'A population statistic about preschool student role.'; D = 'student or trainee role population statistic' },
    @{ Row = 737; C = 'A population statistic about protective factor for harmful behaviour.'; D = 'disposition population statistic' },
    @{ Row = 744; C = 'A population statistic about quantity of valuable material resource owned.'; D = 'data item population statistic' },
    @{ Row = 749; C = 'A population statistic about queer.'; D = 'sexual orientation population statistic' },
    @{ Row = 752; C = 'A population statistic about questioning sexual orientation.'; D = 'self-identity population statistic' },
    @{ Row = 755; C = 'A population statistic about relationship status.'; D = 'personal attribute population statistic' },
    @{ Row = 758; C = 'A population statistic about religious group membership.'; D = 'personal attribute population statistic' },
    @{ Row = 761; C = 'A population statistic about rent-free occupier.'; D = 'person population statistic' },
    @{ Row = 764; C = 'A population statistic about rent-free occupier without owner''s permission.'; D = 'rent-free occupier population statistic' },
    @{ Row = 767; C = 'A population statistic about renter.'; D = 'person population statistic' },
    @{ Row = 770; C = 'A population statistic about renter of housing from a social provider.'; D = 'renter population statistic' },
    @{ Row = 773; C = 'A population statistic about residential facility owner.'; D = 'owner population statistic' },
    @{ Row = 776; C = 'A population statistic about retired status.'; D = 'personal attribute population statistic' },
    @{ Row = 779; C = 'A population statistic about school student role.'; D = 'student or trainee role population statistic' },
    @{ Row = 782; C = 'A population statistic about second generation immigrant.'; D = 'person population statistic' },
    @{ Row = 785; C = 'A population statistic about self employed status.'; D = 'employment status population statistic' },
    @{ Row = 788; C = 'A population statistic about sexual orientation.'; D = 'personal attribute population statistic' },
    @{ Row = 791; C = 'A population statistic about sibling.'; D = 'family member population statistic' },
    @{ Row = 794; C = 'A population statistic about single.'; D = 'relationship status population statistic' },
    @{ Row = 797; C = 'A population statistic about socioeconomic status category.'; D = 'data item population statistic' },
    @{ Row = 800; C = 'A population statistic about socioeconomic status score.'; D = 'data item population statistic' },
    @{ Row = 805; C = 'A population statistic about stay at home parent or guardian status.'; D = 'personal attribute population statistic' },
    @{ Row = 808; C = 'A population statistic about step-parent.'; D = 'parent population statistic' },
    @{ Row = 811; C = 'A population statistic about step-sibling.'; D = 'sibling population statistic' },
    @{ Row = 814; C = 'A population statistic about stepbrother.'; D = 'step-sibling population statistic' },
    @{ Row = 817; C = 'A population statistic about stepchild.'; D = 'child relation population statistic' },
    @{ Row = 820; C = 'A population statistic about stepdaughter.'; D = 'stepchild population statistic' },
    @{ Row = 823; C = 'A population statistic about stepfather.'; D = 'step-parent population statistic' },
    @{ Row = 826; C = 'A population statistic about stepmother.'; D = 'step-parent population statistic' },
    @{ Row = 829; C = 'A population statistic about stepsister.'; D = 'step-sibling population statistic' },
    @{ Row = 832; C = 'A population statistic about stepson.'; D = 'stepchild population statistic' },
    @{ Row = 835; C = 'A population statistic about student or trainee role.'; D = 'role population statistic' },
    @{ Row = 838; C = 'A population statistic about teenager.'; D = 'person population statistic' },
    @{ Row = 841; C = 'A population statistic about transgender.'; D = 'gender identity population statistic' },
    @{ Row = 844; C = 'A population statistic about twin.'; D = 'sibling population statistic' },
    @{ Row = 847; C = 'A population statistic about unawareness of a behaviour.'; D = 'situational personal attribute population statistic' },
    @{ Row = 850; C = 'A population statistic about uncle.'; D = 'family member population statistic' },
    @{ Row = 853; C = 'A population statistic about undecidedness about enacting a behaviour.'; D = 'situational personal attribute population statistic' },
    @{ Row = 856; C = 'A population statistic about undergraduate student role.'; D = 'higher education student role population statistic' },
    @{ Row = 859; C = 'A population statistic about unemployed.'; D = 'employment status population statistic' },
    @{ Row = 862; C = 'A population statistic about unpaid carer for an adult status.'; D = 'personal attribute population statistic' },
    @{ Row = 865; C = 'A population statistic about value of valuable material resource owned.'; D = 'data item population statistic' },
    @{ Row = 870; C = 'A population statistic about vocational training student or trainee role.'; D = 'student or trainee role population statistic' },
    @{ Row = 873; C = 'A population statistic about voluntary worker status.'; D = 'personal attribute population statistic' },
    @{ Row = 876; C = 'A population statistic about widowed.'; D = 'relationship status population statistic' },
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
}
